$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header in D1 from "correctANS" to "correctresponse"
$ws.Range("D1").Value = "correctresponse"

# Update the active selection to D1 (matches recorded sheetView selection)
$ws.Range("D1").Select()
